# Nowcast update: "add results from latest run".
#
# Appends the newest nowcast run (row for 2025-08-30) and refreshes the
# Prognose / Revision figures for every existing date row, matching the
# latest model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged labels, rewritten defensively) -----------------
$headers = @(
    "Row", "Prognose", "surveys", "production", "orders", "turnover",
    "financial", "labor market", "prices", "national accounts", "Revision"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Data rows --------------------------------------------------------
# Each entry: row label (date) followed by the 10 numeric columns
# (Prognose, surveys, production, orders, turnover, financial,
#  labor market, prices, national accounts, Revision).
$rows = @(
    ,@("2025-03-30", "0.27190224152723685", "0", "0", "0", "0", "0", "0", "0", "0", "0")
    ,@("2025-04-15", "0.27762423439395828", "0", "-0.0050352452070996723", "0.00010123641654144588", "-0.0011565721366977259", "0.00010312515053895131", "-6.7135260714282479e-05", "-0.00022137089671385725", "0", "-0.0008580611835595775")
    ,@("2025-04-30", "0.27312898546371989", "-0.0019015034224111596", "0", "0.00055176303277792982", "9.9907966260616351e-05", "0", "0.0001209920247802142", "-0.0030871138455360911", "0.00035352187958286252", "-0.00052185962890566318")
    ,@("2025-05-15", "0.28675355805503216", "0.01079421143679549", "-0.0070197991412063966", "0.00034976274163227175", "0.0012089853445363477", "-0.0015260485608578731", "5.3182464915103453e-06", "-0.00069297105011207802", "0", "-0.0007542776137148266")
    ,@("2025-05-30", "0.21915614346304746", "0.021718509437570633", "0", "-0.00042743729192465331", "2.3231732412905354e-05", "0", "-0.00010820641071184305", "-0.0025634545427196346", "0", "-0.0064174365488265939")
    ,@("2025-06-15", "0.222553466560727", "0", "-0.0027435837580150921", "-0.0013753478828626474", "-0.0073033289792951575", "0.0013371127670396639", "0", "0.00039198436765868367", "0", "-0.001309619237427645")
    ,@("2025-06-30", "0.26781604447333335", "-0.054981178824773866", "0", "9.6933297367771657e-05", "-0.00014600928447466548", "0", "5.7744465132277192e-05", "0.00058430443242319936", "0", "0.0011697722660345133")
    ,@("2025-07-15", "0.27542285401373273", "0", "0.006431796016235369", "-0.0032915808941957267", "-0.0042442079923689984", "0.0023494281690730843", "-0.00017751350632216512", "0.00023351218478545542", "0", "-0.00023935840994854241")
    ,@("2025-07-30", "0.23088365020532678", "0.094186412990487745", "0", "-0.00091030819026355673", "-0.00010897130556109803", "0", "-1.6650701556673681e-05", "0.0011234820384054506", "-0.0025089372037770748", "0.0060865566002865079")
    ,@("2025-08-15", "0.27393717556545016", "0", "-0.03407729415899588", "0.0034790442955649576", "0.0056212747207544014", "0.0040413386042110755", "0.00091243400586765028", "0.005552675991624849", "0", "0.0086237788170338003")
    ,@("2025-08-30", "0.24229227434531142", "-0.067390310106783619", "0", "0.0023751904307887351", "3.7415191255075843e-05", "0", "-4.7712224325675936e-06", "-0.0033315497051391884", "0", "-0.014794525487071941")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $vals = $rows[$i]

    # Column A holds a date-formatted label ("2025-03-30"); prefix with an
    # apostrophe so Excel stores it as plain text instead of auto-converting
    # it to a date serial number (matches the source data's string type).
    $ws.Cells.Item($r, 1).Value = "'" + $vals[0]

    for ($c = 1; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$vals[$c]
    }
}
